# The paragraph under "Design Decisions" originally read:
#   "I'm using spring boots. The first reason I'm using the spring boots is
#    because it's mandatory to use spring boots in this semester, but
#    there's other reason why I'm using spring boots. Because on this
#    semester we learn about Dependency Injection from SOLID, Spring
#    boots help us with ..."
# It is trimmed down to:
#   "I'm using spring boots. On this semester we learn about Dependency
#    Injection from SOLID, Spring boots help us with ..."
#
# Build the strings with explicit char codes for the curly apostrophe
# (U+2019) so this is robust regardless of source-file encoding.

$apos = [char]0x2019

$old = "I" + $apos + "m using spring boots. The first reason I" + $apos + `
       "m using the spring boots is because it" + $apos + "s mandatory " + `
       "to use spring boots in this semester, but there" + $apos + "s " + `
       "other reason why I" + $apos + "m using spring boots. Because on " + `
       "this semester we learn about Dependency Injection from SOLID, " + `
       "Spring boots help us with"

$new = "I" + $apos + "m using spring boots. On this semester we learn " + `
       "about Dependency Injection from SOLID, Spring boots help us with"

$d = $word.ActiveDocument

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                              $true, 1, $false, $new, 2)

if (-not $found) {
    throw "edit.ps1: target sentence was not found in the document"
}
